$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new student record - Priyanka Bhardwaj
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "T"
$ws.Range("C3").Value = 76
$ws.Range("D3").Value = "Priyanka Bhardwaj"
$ws.Range("E3").Value = "anubhav.patrick@kiet.edu"
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:anubhav.patrick@kiet.edu")
$ws.Rows.Item(3).RowHeight = 15.75

# Row 4: new student record - Tishka Gupra
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "J"
$ws.Range("C4").Value = 88
$ws.Range("D4").Value = "Tishka Gupra"
$ws.Range("E4").Value = "anubhavpatrick@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:anubhavpatrick@gmail.com")
$ws.Rows.Item(4).RowHeight = 15.75

# Move the active selection to E5, matching the saved worksheet view state
$ws.Range("E5").Select()
